$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "eta" related headers/values to "Y" naming scheme
$ws.Range("D1").Value = "Y_min"
$ws.Range("E1").Value = "Y_max"
$ws.Range("F1").Value = "Y"

$ws.Range("N2:N15").Value = "W_asym"
$ws.Range("O2:O15").Value = "Y"

# Update the active selection to match the authored state
$ws.Range("L21").Select()
